$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 175, shifting existing rows 175:252 down to 176:253
$ws.Rows("175").Insert()

# Copy the (now shifted) old row 175 data, which is now at row 176, into new row 175
$ws.Range("A176:R176").Copy() | Out-Null
$ws.Range("A175:R175").PasteSpecial() | Out-Null

# Overwrite the new data point's changed fields in row 175
$ws.Range("D175").Value = "2022-11-10"
$ws.Range("J175").Value = 49
$ws.Range("K175").Value = 16000
$ws.Range("L175").Value = 17000
$ws.Range("M175").Value = 16449
$ws.Range("P175").Value = 5483
